# simple_import_term_missing_lang.xlsx fix:
# Column C1's header was incorrectly duplicating D1's "altLabel" label.
# Introduce a distinct "altLabel_fi" header for column C, leaving D1 (and
# every other header/cell) untouched, and move the active selection from
# D2 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "altLabel_fi"

$ws.Range("C2").Select()
